$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly driver report refresh: updated Critical Minutes and Good Roaming
# Calculation (%) for the "Bad Drivers" entry, and filled in the missing
# Driver Vintage date for the third "Good Drivers" entry.

# Critical Minutes (row 3 detail + row 4 total)
$ws.Range("C3").Value = 294
$ws.Range("C4").Value = 294

# Good Roaming Calculation (%)
$ws.Range("D3").Value = 98.9

# Driver Vintage for Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1.
# Force text formatting first so the date-like string isn't auto-converted
# into a date serial number (matching how the sibling vintage cells, e.g.
# E13/E14, already store their dates as plain text).
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2022-08-29"

# Restore the cell's original General/right-aligned style (shared with
# D12/D13/D14 and E13/E14) now that the text value is locked in, instead of
# keeping the new Text-format style that NumberFormat = "@" minted.
$ws.Range("D12").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
